$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 103, shifting existing rows 103:112 down to 104:113.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new weekly price record.
$ws.Range("A103").Value = 8
$ws.Range("B103").Value = "Terminal La Palmera de La Serena"
$ws.Range("C103").Value = "Coquimbo"
$ws.Range("D103").Value = 44578
$ws.Range("E103").Value = 4
$ws.Range("F103").Value = 100112001
$ws.Range("G103").Value = "Berenjena"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 700
$ws.Range("K103").Value = 8000
$ws.Range("L103").Value = 9000
$ws.Range("M103").Value = 8500
$ws.Range("N103").Value = "$/caja 50 unidades"
$ws.Range("O103").Value = "Región de Arica y Parinacota"
$ws.Range("P103").Value = 170
$ws.Range("Q103").Value = 50
$ws.Range("R103").Value = "Hortaliza"
